$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 742.7857
$ws.Range("J17").Value = 742.7857
$ws.Range("L17").Value = 2228.3571
$ws.Range("N17").Value = -2564.3571

$ws.Range("H106").Value = 1534.4231
$ws.Range("I106").Value = 1495.8
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 1495.8
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -864.8
$ws.Range("N106").Value = -3762

$ws.Range("H112").Value = 2591.5625
$ws.Range("J112").Value = 2825.1724
$ws.Range("L112").Value = 8475.5172
$ws.Range("N112").Value = -10691.5172

$ws.Range("H116").Value = 5077.276
$ws.Range("J116").Value = 4630.125
$ws.Range("L116").Value = 4630.125
$ws.Range("N116").Value = -11514.125

$ws.Range("H118").Value = 63001890
$ws.Range("I118").Value = 168000670
$ws.Range("J118").Value = 2617.6
$ws.Range("K118").Value = 504002010
$ws.Range("L118").Value = 7852.799999999999
$ws.Range("M118").Value = -504000353
$ws.Range("N118").Value = -11166.8

$ws.Range("H135").Value = 33997.387
$ws.Range("I135").Value = 45012.914
$ws.Range("J135").Value = 2327.75
$ws.Range("K135").Value = 405116.226
$ws.Range("L135").Value = 20949.75
$ws.Range("M135").Value = -402581.226
$ws.Range("N135").Value = -26019.75

$ws.Range("H137").Value = 2885702.5
$ws.Range("I137").Value = 1316874.5
$ws.Range("J137").Value = 7143950
$ws.Range("K137").Value = 3950623.5
$ws.Range("L137").Value = 21431850
$ws.Range("M137").Value = -3948073.5
$ws.Range("N137").Value = -21436950

$ws.Range("H138").Value = 1711.7885
$ws.Range("I138").Value = 1079.3684
$ws.Range("J138").Value = 3428.3572
$ws.Range("K138").Value = 3238.1052
$ws.Range("L138").Value = 10285.0716
$ws.Range("M138").Value = 1901.8948
$ws.Range("N138").Value = -20565.0716

$ws.Range("H141").Value = 1681.5254
$ws.Range("I141").Value = 1225
$ws.Range("J141").Value = 2642.6316
$ws.Range("K141").Value = 3675
$ws.Range("L141").Value = 7927.8948
$ws.Range("M141").Value = 1505
$ws.Range("N141").Value = -18287.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3282.24
$ws.Range("I32").Value = 2793.0557
$ws.Range("J32").Value = 7684.9
$ws.Range("K32").Value = 2793.0557
$ws.Range("L32").Value = 7684.9
$ws.Range("M32").Value = -2506.0557
$ws.Range("N32").Value = -8258.9

$ws.Range("H61").Value = 1423.8948
$ws.Range("I61").Value = 1488.2069
$ws.Range("J61").Value = 1216.6666
$ws.Range("K61").Value = 1488.2069
$ws.Range("L61").Value = 1216.6666
$ws.Range("M61").Value = -1276.2069
$ws.Range("N61").Value = -1640.6666

$ws.Range("H74").Value = 963.67926
$ws.Range("I74").Value = 933.5454999999999
$ws.Range("J74").Value = 1111
$ws.Range("K74").Value = 933.5454999999999
$ws.Range("L74").Value = 1111
$ws.Range("M74").Value = -59.54549999999995
$ws.Range("N74").Value = -2859

$ws.Range("H77").Value = 963.67926
$ws.Range("I77").Value = 933.5454999999999
$ws.Range("J77").Value = 1111
$ws.Range("K77").Value = 4667.7275
$ws.Range("L77").Value = 5555
$ws.Range("M77").Value = -299.7275
$ws.Range("N77").Value = -14291

$ws.Range("H111").Value = 25000
$ws.Range("J111").Value = 25000
$ws.Range("L111").Value = 25000
$ws.Range("N111").Value = -33180

$ws.Range("H132").Value = 96175.89
$ws.Range("I132").Value = 101186.44
$ws.Range("K132").Value = 303559.32
$ws.Range("M132").Value = -301029.32

$ws.Range("H136").Value = 1423.8948
$ws.Range("I136").Value = 1488.2069
$ws.Range("J136").Value = 1216.6666
$ws.Range("K136").Value = 4464.620699999999
$ws.Range("L136").Value = 3649.9998
$ws.Range("M136").Value = -1914.620699999999
$ws.Range("N136").Value = -8749.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 474.2353
$ws.Range("I80").Value = 89.666664
$ws.Range("J80").Value = 556.6429000000001
$ws.Range("K80").Value = 89.666664
$ws.Range("L80").Value = 556.6429000000001
$ws.Range("M80").Value = 908.333336
$ws.Range("N80").Value = -2552.6429

$ws.Range("H83").Value = 474.2353
$ws.Range("I83").Value = 89.666664
$ws.Range("J83").Value = 556.6429000000001
$ws.Range("K83").Value = 448.33332
$ws.Range("L83").Value = 2783.2145
$ws.Range("M83").Value = 4543.66668
$ws.Range("N83").Value = -12767.2145

$ws.Range("H109").Value = 19083.334
$ws.Range("J109").Value = 19083.334
$ws.Range("L109").Value = 19083.334
$ws.Range("N109").Value = -21857.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H31").Value = 1822.3489
$ws.Range("I31").Value = 1327.7142
$ws.Range("J31").Value = 3986.375
$ws.Range("K31").Value = 1327.7142
$ws.Range("L31").Value = 3986.375
$ws.Range("M31").Value = -1032.7142
$ws.Range("N31").Value = -4576.375

$ws.Range("H34").Value = 1822.3489
$ws.Range("I34").Value = 1327.7142
$ws.Range("J34").Value = 3986.375
$ws.Range("K34").Value = 1327.7142
$ws.Range("L34").Value = 3986.375
$ws.Range("M34").Value = -1125.7142
$ws.Range("N34").Value = -4390.375

$ws.Range("H58").Value = 1702.079
$ws.Range("I58").Value = 1792.0294
$ws.Range("K58").Value = 1792.0294
$ws.Range("M58").Value = -1589.0294

$ws.Range("H105").Value = 3234.524
$ws.Range("I105").Value = 3800.1
$ws.Range("J105").Value = 2720.3635
$ws.Range("K105").Value = 3800.1
$ws.Range("L105").Value = 2720.3635
$ws.Range("M105").Value = -2053.1
$ws.Range("N105").Value = -6214.363499999999

$ws.Range("H132").Value = 2990.8462
$ws.Range("I132").Value = 2510.9666
$ws.Range("J132").Value = 4590.4443
$ws.Range("K132").Value = 7532.899800000001
$ws.Range("L132").Value = 13771.3329
$ws.Range("M132").Value = -5002.899800000001
$ws.Range("N132").Value = -18831.3329

$ws.Range("H134").Value = 7357.7
$ws.Range("I134").Value = 8200.68
$ws.Range("J134").Value = 3142.8
$ws.Range("K134").Value = 24602.04
$ws.Range("L134").Value = 9428.400000000001
$ws.Range("M134").Value = -22067.04
$ws.Range("N134").Value = -14498.4

$ws.Range("H136").Value = 1702.079
$ws.Range("I136").Value = 1792.0294
$ws.Range("K136").Value = 5376.0882
$ws.Range("M136").Value = -2826.0882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 49802.895
$ws.Range("J122").Value = 58617.74
$ws.Range("L122").Value = 527559.66
$ws.Range("N122").Value = -532459.66

$ws.Range("H131").Value = 938.4400000000001
$ws.Range("J131").Value = 968.64514
$ws.Range("L131").Value = 2905.93542
$ws.Range("N131").Value = -12985.93542

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 60892.31
$ws.Range("J135").Value = 60892.31
$ws.Range("L135").Value = 60892.31
$ws.Range("N135").Value = -71032.31

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 37391.5
$ws.Range("J92").Value = 37391.5
$ws.Range("L92").Value = 37391.5
$ws.Range("N92").Value = -42383.5

$ws.Range("H132").Value = 1437.875
$ws.Range("I132").Value = 1344
$ws.Range("K132").Value = 4032
$ws.Range("M132").Value = -1502

$ws.Range("H136").Value = 1253.2333
$ws.Range("I136").Value = 1166.5491
$ws.Range("J136").Value = 1744.4445
$ws.Range("K136").Value = 3499.6473
$ws.Range("L136").Value = 5233.333500000001
$ws.Range("M136").Value = -949.6472999999996
$ws.Range("N136").Value = -10333.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2506.4333
$ws.Range("I126").Value = 2099.889
$ws.Range("J126").Value = 3116.25
$ws.Range("K126").Value = 6299.667
$ws.Range("L126").Value = 9348.75
$ws.Range("M126").Value = -3829.667
$ws.Range("N126").Value = -14288.75

$ws.Range("H136").Value = 1696.0209
$ws.Range("I136").Value = 1845.027
$ws.Range("J136").Value = 1194.8182
$ws.Range("K136").Value = 5535.081
$ws.Range("L136").Value = 3584.4546
$ws.Range("M136").Value = -2985.081
$ws.Range("N136").Value = -8684.454600000001

$ws.Range("H138").Value = 49249
$ws.Range("J138").Value = 49249
$ws.Range("L138").Value = 49249
$ws.Range("N138").Value = -59529
